$d = $word.ActiveDocument

# Locate the paragraph containing "This is a regular line" and insert two
# new empty (Normal-styled) paragraphs directly after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*This is a regular line*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
    $target.Range.InsertParagraphAfter()
}
